# Updates the crypto price/volume table to the latest scraped values.
# Column D (Price) cells are stored as literal TEXT in the source data
# (e.g. "67.089.78", "0.582"), so every Price write below is given a
# leading apostrophe to force a text literal, exactly like a user typing
# 0.581 into a General-formatted cell, instead of letting Excel
# autodetection silently reinterpret it as a Number/scientific value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = ("'66.916.22")
$ws.Range("E2").Value = ("  +1.95%  ")

$ws.Range("D3").Value = ("'3.277.09")
$ws.Range("E3").Value = ("  -0.93%  ")

$ws.Range("E4").Value = ("  +0.10%  ")

$ws.Range("D5").Value = ("'569.80")
$ws.Range("E5").Value = ("  -1.60%  ")

$ws.Range("D6").Value = ("'175.07")
$ws.Range("E6").Value = ("  -4.88%  ")

$ws.Range("E7").Value = ("  +0.11%  ")

$ws.Range("D8").Value = ("'0.581")

$ws.Range("D9").Value = ("'3.274.37")
$ws.Range("E9").Value = ("  -0.89%  ")

$ws.Range("D10").Value = ("'0.173")
$ws.Range("E10").Value = ("  -2.42%  ")

$ws.Range("D11").Value = ("'0.571")
$ws.Range("E11").Value = ("  -0.16%  ")

$ws.Range("D12").Value = ("'45.58")
$ws.Range("E12").Value = ("  -2.30%  ")

$ws.Range("D13").Value = ("'0.0000268")
$ws.Range("E13").Value = ("  +1.10%  ")

$ws.Range("D14").Value = ("'687.54")
$ws.Range("E14").Value = ("  +8.43%  ")

$ws.Range("D15").Value = ("'3.809.57")
$ws.Range("E15").Value = ("  -0.71%  ")

$ws.Range("D16").Value = ("'8.29")
$ws.Range("E16").Value = ("  -1.95%  ")

$ws.Range("D17").Value = ("'67.038.52")
$ws.Range("E17").Value = ("  +1.82%  ")

$ws.Range("D18").Value = ("'0.118")
$ws.Range("E18").Value = ("  +0.93%  ")

$ws.Range("D19").Value = ("'3.286.66")
$ws.Range("E19").Value = ("  -0.63%  ")

$ws.Range("D20").Value = ("'17.29")
$ws.Range("E20").Value = ("  -3.32%  ")

$ws.Range("D21").Value = ("'10.70")
$ws.Range("E21").Value = ("  -2.76%  ")

$ws.Range("D22").Value = ("'0.886")
$ws.Range("E22").Value = ("  -0.55%  ")

$ws.Range("D23").Value = ("'16.90")
$ws.Range("E23").Value = ("  -4.56%  ")

$ws.Range("D24").Value = ("'5.10")
$ws.Range("E24").Value = ("  +1.75%  ")

$ws.Range("D25").Value = ("'98.93")
$ws.Range("E25").Value = ("  -1.42%  ")

$ws.Range("D26").Value = ("'3.87")
$ws.Range("E26").Value = ("  -2.40%  ")

$ws.Range("D27").Value = ("'2.69")
$ws.Range("E27").Value = ("  -1.23%  ")

$ws.Range("D28").Value = ("'9.26")
$ws.Range("E28").Value = ("  -1.51%  ")

$ws.Range("D29").Value = ("'32.87")
$ws.Range("E29").Value = ("  +6.59%  ")

$ws.Range("D30").Value = ("'8.35")
$ws.Range("E30").Value = ("  +0.04%  ")

$ws.Range("D31").Value = ("'6.71")
$ws.Range("E31").Value = ("  +2.09%  ")

$ws.Range("D32").Value = ("'573.31")
$ws.Range("E32").Value = ("  -3.59%  ")

$ws.Range("D33").Value = ("'3.875.85")
$ws.Range("E33").Value = ("  +0.58%  ")

$ws.Range("D34").Value = ("'10.80")
$ws.Range("E34").Value = ("  -1.29%  ")

$ws.Range("E35").Value = ("  -2.12%  ")

$ws.Range("D36").Value = ("'0.999")
$ws.Range("E36").Value = ("  -0.17%  ")

$ws.Range("D37").Value = ("'55.52")
$ws.Range("E37").Value = ("  -0.31%  ")

$ws.Range("D38").Value = ("'3.31")
$ws.Range("E38").Value = ("  -11.61%  ")

$ws.Range("E39").Value = ("  +1.50%  ")

$ws.Range("D40").Value = ("'2.60")
$ws.Range("E40").Value = ("  +0.00%  ")

$ws.Range("E41").Value = ("  -1.30%  ")

$ws.Range("D42").Value = ("'31.75")
$ws.Range("E42").Value = ("  -1.74%  ")

$ws.Range("D43").Value = ("'0.0" + [char]0x2083 + "0668")
$ws.Range("E43").Value = ("  -4.41%  ")

$ws.Range("D44").Value = ("'3.00")
$ws.Range("E44").Value = ("  -4.36%  ")

$ws.Range("D45").Value = ("'0.326")
$ws.Range("E45").Value = ("  -1.83%  ")

$ws.Range("D46").Value = ("'0.0404")
$ws.Range("E46").Value = ("  -0.94%  ")

$ws.Range("E47").Value = ("  +0.19%  ")

$ws.Range("B49").Value = ("Mantle")
$ws.Range("C49").Value = ("https://coinranking.com/coin/BoI4ux0nd+mantle-mnt")
$ws.Range("D49").Value = ("'1.37")
$ws.Range("E49").Value = ("  +7.20%  ")

$ws.Range("B50").Value = ("ThetaToken")
$ws.Range("C50").Value = ("https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta")
$ws.Range("D50").Value = ("'2.52")
$ws.Range("E50").Value = ("  -0.27%  ")

$ws.Range("D51").Value = ("'130.44")
$ws.Range("E51").Value = ("  -0.28%  ")
